$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format to preserve numeric-looking strings (e.g. "1.001", "0.4021")
# exactly as text instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.630.28"
$ws.Range("E2").Value = "  +3.60%  "

$ws.Range("D3").Value = "1.696.69"
$ws.Range("E3").Value = "  +2.01%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "318.41"
$ws.Range("E5").Value = "  +2.72%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("E7").Value = "  +1.55%  "

$ws.Range("D8").Value = "0.4021"
$ws.Range("E8").Value = "  +1.47%  "

$ws.Range("D9").Value = "1.539"
$ws.Range("E9").Value = "  +8.13%  "

$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("D11").Value = "53.17"
$ws.Range("E11").Value = "  +7.25%  "

$ws.Range("D12").Value = "0.08783"
$ws.Range("E12").Value = "  +1.35%  "

$ws.Range("D13").Value = "7.275"
$ws.Range("E13").Value = "  +8.76%  "

$ws.Range("D14").Value = "23.28"
$ws.Range("E14").Value = "  +2.36%  "

$ws.Range("D15").Value = "0.00001321"
$ws.Range("E15").Value = "  +0.47%  "

$ws.Range("D16").Value = "7.630"
$ws.Range("E16").Value = "  +5.24%  "

$ws.Range("D17").Value = "1.694.70"
$ws.Range("E17").Value = "  +1.95%  "

$ws.Range("D18").Value = "101.09"
$ws.Range("E18").Value = "  +1.00%  "

$ws.Range("D19").Value = "0.07013"
$ws.Range("E19").Value = "  +3.51%  "

$ws.Range("D20").Value = "19.73"
$ws.Range("E20").Value = "  +3.15%  "

$ws.Range("D21").Value = "6.910"
$ws.Range("E21").Value = "  +3.64%  "

$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("E23").Value = "  +1.69%  "

$ws.Range("D24").Value = "24.623.29"
$ws.Range("E24").Value = "  +3.68%  "

$ws.Range("D25").Value = "3.075"
$ws.Range("E25").Value = "  +8.79%  "

$ws.Range("D26").Value = "2.342"
$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("D28").Value = "159.74"
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").Value = "5.236"
$ws.Range("E29").Value = "  +1.38%  "

$ws.Range("D30").Value = "134.65"
$ws.Range("E30").Value = "  +3.78%  "

$ws.Range("D31").Value = "7.485"
$ws.Range("E31").Value = "  +15.68%  "

$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "1.882.64"
$ws.Range("E32").Value = "  +1.80%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "1.105"
$ws.Range("E33").Value = "  -2.92%  "

$ws.Range("D34").Value = "7.490"
$ws.Range("E34").Value = "  +15.17%  "

$ws.Range("D35").Value = "0.08539"
$ws.Range("E35").Value = "  -0.90%  "

$ws.Range("D36").Value = "11.50"
$ws.Range("E36").Value = "  +11.17%  "

$ws.Range("D37").Value = "1.979"
$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").Value = "0.2742"
$ws.Range("E38").Value = "  +3.13%  "

$ws.Range("D39").Value = "14.58"
$ws.Range("E39").Value = "  +1.09%  "

$ws.Range("D40").Value = "0.02779"
$ws.Range("E40").Value = "  +9.66%  "

$ws.Range("D41").Value = "0.09025"
$ws.Range("E41").Value = "  +2.75%  "

$ws.Range("D42").Value = "1.465"
$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("D43").Value = "0.7728"
$ws.Range("E43").Value = "  +2.16%  "

$ws.Range("D44").Value = "0.7240"
$ws.Range("E44").Value = "  +2.92%  "

$ws.Range("D45").Value = "15.46"
$ws.Range("E45").Value = "  +3.64%  "

$ws.Range("D46").Value = "2.534"
$ws.Range("E46").Value = "  +5.87%  "

$ws.Range("D47").Value = "4.232"
$ws.Range("E47").Value = "  +3.34%  "

$ws.Range("D48").Value = "1.354"
$ws.Range("E48").Value = "  +13.31%  "

$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D50").Value = "141.23"
$ws.Range("E50").Value = "  +1.85%  "

$ws.Range("D51").Value = "0.08047"
$ws.Range("E51").Value = "  +3.61%  "
